# This edit reorders (re-sorts) the 29 data rows of the "Rabanito" sheet
# (rows 2-30) into a new row order. Every underlying record (date, volume,
# prices, origin, quality, etc.) already exists in the sheet; the commit
# simply re-sequences which physical row each record occupies, matching a
# "weekly" re-sort pass (commit message: "Fruta / hortaliza, semanal").
#
# Rather than physically moving rows (which could disturb styles/formats),
# we overwrite each cell that differs between the old row order and the
# new row order with its target value - the net effect on the saved
# workbook is identical to performing the reorder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44231
$ws.Range("J2").Value = 12000
$ws.Range("D3").Value = 44847
$ws.Range("J3").Value = 7900
$ws.Range("D4").Value = 44232
$ws.Range("J4").Value = 16000
$ws.Range("D5").Value = 44230
$ws.Range("J5").Value = 16000
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = 3000
$ws.Range("O5").Value = 'Provincia de Chacabuco'
$ws.Range("P5").Value = 30
$ws.Range("D6").Value = 44846
$ws.Range("J6").Value = 7900
$ws.Range("D7").Value = 44602
$ws.Range("D8").Value = 44602
$ws.Range("I8").Value = 'Segunda'
$ws.Range("J8").Value = 6000
$ws.Range("K8").Value = 2500
$ws.Range("L8").Value = 2500
$ws.Range("M8").Value = 2500
$ws.Range("P8").Value = 25
$ws.Range("D9").Value = 44204
$ws.Range("J9").Value = 7000
$ws.Range("D10").Value = 44188
$ws.Range("J10").Value = 12000
$ws.Range("D11").Value = 44215
$ws.Range("J11").Value = 16000
$ws.Range("O11").Value = 'Provincia de Chacabuco'
$ws.Range("D12").Value = 44162
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 7000
$ws.Range("K12").Value = 3000
$ws.Range("L12").Value = 3000
$ws.Range("M12").Value = 3000
$ws.Range("O12").Value = 'Provincia de Chacabuco'
$ws.Range("P12").Value = 30
$ws.Range("D13").Value = 44245
$ws.Range("J13").Value = 9000
$ws.Range("O13").Value = 'Región Metropolitana'
$ws.Range("D14").Value = 44245
$ws.Range("J14").Value = 5000
$ws.Range("O14").Value = 'Región Metropolitana'
$ws.Range("D15").Value = 44210
$ws.Range("J15").Value = 8800
$ws.Range("K15").Value = 2500
$ws.Range("M15").Value = 2750
$ws.Range("P15").Value = 28
$ws.Range("D16").Value = 44159
$ws.Range("K16").Value = 3000
$ws.Range("M16").Value = 3000
$ws.Range("P16").Value = 30
$ws.Range("D17").Value = 44160
$ws.Range("J17").Value = 7000
$ws.Range("D18").Value = 44161
$ws.Range("J18").Value = 7000
$ws.Range("D19").Value = 44189
$ws.Range("J19").Value = 16000
$ws.Range("D20").Value = 44229
$ws.Range("J20").Value = 16000
$ws.Range("D21").Value = 44186
$ws.Range("J21").Value = 10000
$ws.Range("D22").Value = 44181
$ws.Range("J22").Value = 12000
$ws.Range("D23").Value = 44187
$ws.Range("J23").Value = 12000
$ws.Range("K23").Value = 3000
$ws.Range("M23").Value = 3000
$ws.Range("P23").Value = 30
$ws.Range("D25").Value = 44168
$ws.Range("D26").Value = 44214
$ws.Range("J26").Value = 7000
$ws.Range("D27").Value = 44167
$ws.Range("D28").Value = 44209
$ws.Range("K28").Value = 2500
$ws.Range("M28").Value = 2750
$ws.Range("P28").Value = 28
$ws.Range("D29").Value = 44845
$ws.Range("J29").Value = 7900
$ws.Range("D30").Value = 44600
$ws.Range("J30").Value = 1300
$ws.Range("K30").Value = 3500
$ws.Range("L30").Value = 4000
$ws.Range("M30").Value = 3808
$ws.Range("O30").Value = 'Región Metropolitana'
$ws.Range("P30").Value = 38
